# Updated Arc Reactor UI
# Rebuilds the "ItemsPerMin" recipe table with new rows/values per the
# commit. Operates purely through cell Value/Formula writes plus
# Merge/UnMerge of the small 2x1 "header spans two rows" blocks that this
# sheet uses to show one input/output pair at a time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 4-5 block: Element 65 + Sulfuric Acid -> RefinedE65
# ---------------------------------------------------------------------
$ws.Range("D4").Value = 1
$ws.Range("G4").Value = 1000
$ws.Range("I4").Value = 3

$ws.Range("D5").Value = 1000

# ---------------------------------------------------------------------
# Row 7-8 block: Refined E65 + Uranium -> Unfiltered Deuterium
# ---------------------------------------------------------------------
$ws.Range("D7").Value = 10000
$ws.Range("G7").Value = 10000
$ws.Range("I7").Value = 30

# ---------------------------------------------------------------------
# Row 19-21 block (brand new): LimeStone + Coal + Baux -> Coolant Powder
# ---------------------------------------------------------------------
$ws.Range("C19").Value = "LimeStone"
$ws.Range("D19").Value = 2
$ws.Range("E19").Formula = '=(60/I$19)*D19'
$ws.Range("F19").Value = "Coolant Powder"
$ws.Range("G19").Value = 5
$ws.Range("H19").Formula = '=(60/I$19)*G19'
$ws.Range("I19").Value = 3

$ws.Range("C20").Value = "Coal"
$ws.Range("D20").Value = 1
$ws.Range("E20").Formula = '=(60/I$19)*D20'

$ws.Range("C21").Value = "Baux"
$ws.Range("D21").Value = 2
$ws.Range("E21").Formula = '=(60/I$19)*D21'

$ws.Range("F19:F21").Merge()
$ws.Range("G19:G21").Merge()
$ws.Range("H19:H21").Merge()
$ws.Range("I19:I21").Merge()

# ---------------------------------------------------------------------
# Row 23-24 block (brand new): Water + Coolant Powder -> Cloudy Coolant
# ---------------------------------------------------------------------
$ws.Range("C23").Value = "Water"
$ws.Range("D23").Value = 1000
$ws.Range("E23").Formula = '=(60/I$23)*D23'
$ws.Range("F23").Value = "Cloudy Coolant"
$ws.Range("G23").Value = 1000
$ws.Range("H23").Formula = '=(60/I$23)*G23'
$ws.Range("I23").Value = 12

$ws.Range("C24").Value = "Coolant Powder"
$ws.Range("D24").Value = 20
$ws.Range("E24").Formula = '=(60/I$23)*D24'

$ws.Range("F23:F24").Merge()
$ws.Range("G23:G24").Merge()
$ws.Range("H23:H24").Merge()
$ws.Range("I23:I24").Merge()

# ---------------------------------------------------------------------
# Row 26 (brand new, single-row block): Cloudy Coolant -> Reactor Coolant
# ---------------------------------------------------------------------
$ws.Range("C26").Value = "Cloudy Coolant"
$ws.Range("D26").Value = 1000
$ws.Range("E26").Formula = '=(60/I$26)*D26'
$ws.Range("F26").Value = "Reactor Coolant"
$ws.Range("G26").Value = 1000
$ws.Range("H26").Formula = '=(60/I$26)*G26'
$ws.Range("I26").Value = 24

# ---------------------------------------------------------------------
# Row 10-11 block: Unfiltered Deuterium + Packaged Water -> Liquid
# Deuterium, now gains a second output (Empty Canister) on row 11.
# ---------------------------------------------------------------------
$ws.Range("D10").Value = 1000
$ws.Range("I10").Value = 3

$ws.Range("D11").Value = 2

$ws.Range("F10:F11").UnMerge()
$ws.Range("G10:G11").UnMerge()
$ws.Range("H10:H11").UnMerge()

$ws.Range("F11").Value = "Empty Canister"
$ws.Range("G11").Value = 2
$ws.Range("H11").Formula = '=(60/I$10)*G11'

# ---------------------------------------------------------------------
# Row 13-14 block: now Quartz Crystals + Heavy Modular Frame -> Adv
# Reactor Casing (previously Liquid Deuterium + Adv Reactor Casing ->
# Adv Reactor Core).
# ---------------------------------------------------------------------
$ws.Range("C13").Value = "Quartz Crystals"
$ws.Range("D13").Value = 10
$ws.Range("F13").Value = "Adv Reactor Casing"
$ws.Range("G13").Value = 1
$ws.Range("I13").Value = 6

$ws.Range("C14").Value = "Heavy Modular Frame"
$ws.Range("D14").Value = 1

# ---------------------------------------------------------------------
# Row 16-17 block: now Liquid Deuterium + Adv Reactor Casing -> Adv
# Reactor Core (previously Water -> Reactor Coolant, single input row).
# Row 17 gains a second input (Adv Reactor Casing).
# ---------------------------------------------------------------------
$ws.Range("C16").Value = "Liquid Deuterium"
$ws.Range("D16").Value = 2000
$ws.Range("F16").Value = "Adv Reactor Core"
$ws.Range("G16").Value = 1
$ws.Range("I16").Value = 12

$ws.Range("C17").Value = "Adv Reactor Casing"
$ws.Range("D17").Value = 1
$ws.Range("E17").Formula = '=(60/I$16)*D17'

# ---------------------------------------------------------------------
# Extend the used range down to row 36 with two blank spacer rows, same
# as the rest of the sheet's blank rows.
# ---------------------------------------------------------------------
$ws.Range("C35:I36").HorizontalAlignment = -4108
$ws.Range("C35:I36").VerticalAlignment = -4108

# Selection marker left where the author's cursor ended up.
$ws.Range("F11").Select()
